# Excel template bug fix
#
# 1. On sheet "Bico", clear the "Validado com sucesso! ..." message from
#    column H (Obs_relatorio) for rows 2-8, since these bicos had a real
#    divergence flagged instead (handled by sheet "Tanque").
# 2. On sheet "Tanque", replace the "Validado com sucesso! ..." message in
#    column F (Obs_relatorio) for rows 2-5 with the actual divergence
#    message comparing the SPED value against the report value.

$wb = $excel.ActiveWorkbook

$wsBico = $wb.Worksheets.Item("Bico")
foreach ($r in 2..8) {
    $wsBico.Cells.Item($r, 8).Value = ""
}

$wsTanque = $wb.Worksheets.Item("Tanque")
$wsTanque.Range("F2").Value = "Divergência entre o SPED(7851,80) e o relatório(8793,11)!"
$wsTanque.Range("F3").Value = "Divergência entre o SPED(7851,80) e o relatório(10730,13)!"
$wsTanque.Range("F4").Value = "Divergência entre o SPED(7851,80) e o relatório(10156,99)!"
$wsTanque.Range("F5").Value = "Divergência entre o SPED(7851,80) e o relatório(5418,69)!"
